$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells receiving numeric-looking values
# so Excel keeps them as text (matching the source column type).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "64.266.44"
$ws.Range("E2").Value = "  -2.32%  "

$ws.Range("D3").Value = "3.186.47"
$ws.Range("E3").Value = "  -7.42%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "562.96"
$ws.Range("E5").Value = "  -3.53%  "

$ws.Range("D6").Value = "170.56"
$ws.Range("E6").Value = "  -1.42%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "3.183.47"
$ws.Range("E9").Value = "  -7.44%  "

$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -5.66%  "

$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  -4.11%  "

$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -3.27%  "

$ws.Range("D13").Value = "3.741.35"
$ws.Range("E13").Value = "  -7.31%  "

$ws.Range("E14").Value = "  +0.85%  "

$ws.Range("D15").Value = "27.39"
$ws.Range("E15").Value = "  -4.60%  "

$ws.Range("D16").Value = "64.253.74"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -4.69%  "

$ws.Range("D18").Value = "3.190.68"
$ws.Range("E18").Value = "  -7.29%  "

$ws.Range("D19").Value = "5.67"
$ws.Range("E19").Value = "  -4.26%  "

$ws.Range("D20").Value = "13.10"
$ws.Range("E20").Value = "  -4.77%  "

$ws.Range("D21").Value = "352.19"
$ws.Range("E21").Value = "  -4.14%  "

$ws.Range("E22").Value = "  -5.11%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "69.22"
$ws.Range("E24").Value = "  -3.94%  "

$ws.Range("E25").Value = "  -4.78%  "

$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").Value = "5.65"
$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "1.89"
$ws.Range("E32").Value = "  -4.03%  "

$ws.Range("D33").Value = "22.08"
$ws.Range("E33").Value = "  -6.42%  "

$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("E35").Value = "  -6.37%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "156.07"
$ws.Range("E36").Value = "  -2.86%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.44"
$ws.Range("E37").Value = "  -5.40%  "

$ws.Range("D38").Value = "0.810"
$ws.Range("E38").Value = "  -7.78%  "

$ws.Range("D39").Value = "25.94"
$ws.Range("E39").Value = "  -8.99%  "

$ws.Range("D40").Value = "2.49"
$ws.Range("E40").Value = "  -4.26%  "

$ws.Range("D41").Value = "1.68"
$ws.Range("E41").Value = "  -4.41%  "

$ws.Range("D42").Value = "2.658.26"
$ws.Range("E42").Value = "  -3.49%  "

$ws.Range("D43").Value = "4.15"
$ws.Range("E43").Value = "  -6.45%  "

$ws.Range("D44").Value = "6.00"
$ws.Range("E44").Value = "  -7.13%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0652"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "327.75"
$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("D47").Value = "38.63"
$ws.Range("E47").Value = "  -4.10%  "

$ws.Range("D48").Value = "23.62"
$ws.Range("E48").Value = "  -4.22%  "

$ws.Range("E49").Value = "  -6.01%  "

$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("E51").Value = "  -0.04%  "

Write-Output "applied cryptos update"